# Update the "symbol list" (cryptos.xlsx) values as scraped on
# Sat Dec 17 15:49:06 UTC 2022 by the GitHub Actions job.
#
# The sheet stores every data cell (columns B-G) as literal text, even
# when the text looks like a number (prices, volumes, ...). Writing a
# numeric-looking string straight into Range.Value makes Excel infer a
# Number cell, which both changes the stored cell type and can mangle
# precision (e.g. "237.23" -> 236.94999999999999). To keep these cells
# textual we temporarily mark them as Text (NumberFormat "@") before
# assigning, then restore the original (General / unstyled) look by
# copying the style from the neighboring Link cell in column C, which
# carries the sheet's default, un-styled look for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2";  Value = "236.95" },
    @{ Addr = "D3";  Value = "21.80" },
    @{ Addr = "D4";  Value = "5.384" },
    @{ Addr = "D6";  Value = "6.472" },
    @{ Addr = "D7";  Value = "3.350" },
    @{ Addr = "D8";  Value = "0.7994" },
    @{ Addr = "D9";  Value = "1.037" },
    @{ Addr = "D10"; Value = "0.1388" },
    @{ Addr = "D11"; Value = "0.07304" },
    @{ Addr = "D12"; Value = "0.03126" },
    @{ Addr = "D13"; Value = "0.02973" },
    @{ Addr = "D14"; Value = "0.09240" },
    @{ Addr = "D15"; Value = "0.001664" },
    @{ Addr = "D16"; Value = "3.253" },
    @{ Addr = "D17"; Value = "0.04772" },
    @{ Addr = "D18"; Value = "0.0005714" },
    @{ Addr = "E18"; Value = "17OneONEWorstin24h" },
    @{ Addr = "D19"; Value = "0.006207" },
    @{ Addr = "D20"; Value = "0.005076" },
    @{ Addr = "D21"; Value = "0.001051" },
    @{ Addr = "D22"; Value = "0.0001500" },
    @{ Addr = "D23"; Value = "0.0003903" },
    @{ Addr = "D24"; Value = "3.956" },
    @{ Addr = "D25"; Value = "2.201" },
    @{ Addr = "D27"; Value = "0.1060" },
    @{ Addr = "D40"; Value = "0.04084" },
    @{ Addr = "D41"; Value = "0.006946" },
    @{ Addr = "D42"; Value = "0.003501" },
    @{ Addr = "D43"; Value = "0.1039" },
    @{ Addr = "D44"; Value = "0.008828" },
    @{ Addr = "D45"; Value = "0.00005434" },
    @{ Addr = "D47"; Value = "0.6757" },
    @{ Addr = "D48"; Value = "0.03701" },
    @{ Addr = "E48"; Value = "47BOLOBOLO" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $newValue = $u.Value

    $looksNumeric = $newValue -match '^-?\d+(\.\d+)?$'

    if ($looksNumeric) {
        # Force Excel to keep this as a literal text string instead of
        # inferring a Number cell (which would also round the value).
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        # Restore the row's default (un-styled) appearance, matching the
        # rest of the sheet, by borrowing column C's style for that row.
        $refCell = $ws.Cells.Item($cell.Row, 3)
        $cell.Style = $refCell.Style
    } else {
        $cell.Value = $newValue
    }
}
